# Exchange rate workbook update
# - Sheet "euro" (sheet1): insert a new leading daily NBP mid-rate row
#   (040/A/NBP/2021, 2021-03-01) ahead of the existing 041/A/NBP/2021 row,
#   then append rows for 042, 043 and 044/A/NBP/2021.
# - Sheet "dolar amerykanski" (sheet2): same date/table progression, with
#   this currency's own mid rates.

$wb = $excel.ActiveWorkbook

$xlCenter = -4108

function Set-RateRow {
    param(
        $ws,
        [int]$r,
        [string]$tableNo,
        [string]$effectiveDate,
        [string]$mid
    )

    # D column: reference numbers like "041/A/NBP/2021" are unambiguous text already.
    $ws.Range("D$r").Value = $tableNo
    # E/F columns hold date-like / decimal-like text that Excel would otherwise
    # coerce to a date serial / number, so force literal text with a quote prefix.
    $ws.Range("E$r").Value = "'" + $effectiveDate
    $ws.Range("F$r").Value = "'" + $mid

    $ws.Range("D$r").HorizontalAlignment = $xlCenter
    $ws.Range("D$r").VerticalAlignment = $xlCenter
    $ws.Range("E$r").HorizontalAlignment = $xlCenter
    $ws.Range("E$r").VerticalAlignment = $xlCenter
    $ws.Range("F$r").HorizontalAlignment = $xlCenter
    $ws.Range("F$r").VerticalAlignment = $xlCenter
}

# ---- Sheet 1: euro ----
$ws1 = $wb.Worksheets.Item(1)

# Row 3 now becomes the 040/A/NBP/2021 (2021-03-01) entry ...
Set-RateRow $ws1 3 "040/A/NBP/2021" "2021-03-01" "4.5231"
# ... and the rate table grows with the following days.
Set-RateRow $ws1 4 "041/A/NBP/2021" "2021-03-02" "4.5345"
Set-RateRow $ws1 5 "042/A/NBP/2021" "2021-03-03" "4.5393"
Set-RateRow $ws1 6 "043/A/NBP/2021" "2021-03-04" "4.554"
Set-RateRow $ws1 7 "044/A/NBP/2021" "2021-03-05" "4.5793"

# ---- Sheet 2: dolar amerykanski ----
$ws2 = $wb.Worksheets.Item(2)

Set-RateRow $ws2 3 "040/A/NBP/2021" "2021-03-01" "3.7572"
Set-RateRow $ws2 4 "041/A/NBP/2021" "2021-03-02" "3.7765"
Set-RateRow $ws2 5 "042/A/NBP/2021" "2021-03-03" "3.7509"
